# Scheduled market-price refresh: updates the raw price/profit columns
# (H..N) on each job sheet's Leve-profit table with freshly scraped
# Yojimbo market data. Columns: H currentAveragePrice, I/J NQ/HQ average
# price, K/L NQ/HQ leve price, M/N NQ/HQ leve profit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 244.82608
$ws.Range("I11").Value = 244.82608
$ws.Range("K11").Value = 244.82608
$ws.Range("M11").Value = -104.82608
$ws.Range("H32").Value = 2428.0386
$ws.Range("I32").Value = 2702.4167
$ws.Range("J32").Value = 2192.8572
$ws.Range("K32").Value = 2702.4167
$ws.Range("L32").Value = 2192.8572
$ws.Range("M32").Value = -2376.4167
$ws.Range("N32").Value = -2844.8572
$ws.Range("H40").Value = 999.9231
$ws.Range("J40").Value = 1000
$ws.Range("L40").Value = 1000
$ws.Range("N40").Value = -1350
$ws.Range("H43").Value = 1236
$ws.Range("J43").Value = 1155.4286
$ws.Range("L43").Value = 1155.4286
$ws.Range("N43").Value = -1293.4286
$ws.Range("H51").Value = 5020
$ws.Range("I51").Value = 2300
$ws.Range("J51").Value = 5700
$ws.Range("K51").Value = 2300
$ws.Range("L51").Value = 5700
$ws.Range("M51").Value = -1816
$ws.Range("N51").Value = -6668
$ws.Range("H55").Value = 199.55556
$ws.Range("I55").Value = 82.25
$ws.Range("J55").Value = 293.4
$ws.Range("K55").Value = 82.25
$ws.Range("L55").Value = 293.4
$ws.Range("M55").Value = 131.75
$ws.Range("N55").Value = -721.4
$ws.Range("H121").Value = 770
$ws.Range("I121").Value = 770
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 2310
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -563
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 2505264.8
$ws.Range("I132").Value = 2722550.8
$ws.Range("J132").Value = 6475
$ws.Range("K132").Value = 8167652.399999999
$ws.Range("L132").Value = 19425
$ws.Range("M132").Value = -8165122.399999999
$ws.Range("N132").Value = -24485
$ws.Range("H135").Value = 5396
$ws.Range("I135").Value = 845.7778
$ws.Range("J135").Value = 11696.308
$ws.Range("K135").Value = 7612.000199999999
$ws.Range("L135").Value = 105266.772
$ws.Range("M135").Value = -5077.000199999999
$ws.Range("N135").Value = -110336.772
$ws.Range("H140").Value = 43238.46
$ws.Range("J140").Value = 43238.46
$ws.Range("L140").Value = 43238.46
$ws.Range("N140").Value = -53598.46

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9303.166999999999
$ws.Range("I74").Value = 9874.909
$ws.Range("J74").Value = 3014
$ws.Range("K74").Value = 9874.909
$ws.Range("L74").Value = 3014
$ws.Range("M74").Value = -9000.909
$ws.Range("N74").Value = -4762
$ws.Range("H77").Value = 9303.166999999999
$ws.Range("I77").Value = 9874.909
$ws.Range("J77").Value = 3014
$ws.Range("K77").Value = 49374.545
$ws.Range("L77").Value = 15070
$ws.Range("M77").Value = -45006.545
$ws.Range("N77").Value = -23806
$ws.Range("H97").Value = 2146.111
$ws.Range("I97").Value = 907.7143
$ws.Range("K97").Value = 907.7143
$ws.Range("M97").Value = -411.7143
$ws.Range("H102").Value = 6000
$ws.Range("I102").Value = 6000
$ws.Range("K102").Value = 6000
$ws.Range("M102").Value = -4378
$ws.Range("H110").Value = 938
$ws.Range("I110").Value = 938
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 938
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1107
$ws.Range("N110").ClearContents()
$ws.Range("H133").Value = 475380.5
$ws.Range("J133").Value = 475380.5
$ws.Range("L133").Value = 475380.5
$ws.Range("N133").Value = -480440.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H104").Value = 40684
$ws.Range("J104").Value = 40684
$ws.Range("L104").Value = 40684
$ws.Range("N104").Value = -47672
$ws.Range("H105").Value = 3220.5557
$ws.Range("I105").Value = 3478
$ws.Range("J105").Value = 1933.3334
$ws.Range("K105").Value = 3478
$ws.Range("L105").Value = 1933.3334
$ws.Range("M105").Value = -1731
$ws.Range("N105").Value = -5427.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 30092
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 30092
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 30092
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -31342
$ws.Range("H54").Value = 14980
$ws.Range("J54").Value = 14980
$ws.Range("L54").Value = 14980
$ws.Range("N54").Value = -16296
$ws.Range("H99").Value = 1779.8182
$ws.Range("I99").Value = 1725
$ws.Range("K99").Value = 1725
$ws.Range("M99").Value = -227
$ws.Range("H100").Value = 44780
$ws.Range("J100").Value = 44780
$ws.Range("L100").Value = 44780
$ws.Range("N100").Value = -46944
$ws.Range("H112").Value = 100000
$ws.Range("J112").Value = 100000
$ws.Range("L112").Value = 100000
$ws.Range("N112").Value = -102954
$ws.Range("H126").Value = 1779.8182
$ws.Range("I126").Value = 1725
$ws.Range("K126").Value = 5175
$ws.Range("M126").Value = -2705

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 140.28572
$ws.Range("I12").Value = 130.33333
$ws.Range("J12").Value = 147.75
$ws.Range("K12").Value = 390.99999
$ws.Range("L12").Value = 443.25
$ws.Range("M12").Value = -217.99999
$ws.Range("N12").Value = -789.25
$ws.Range("H113").Value = 778.4643
$ws.Range("I113").Value = 1003.4828
$ws.Range("J113").Value = 536.7778
$ws.Range("K113").Value = 3010.4484
$ws.Range("L113").Value = 1610.3334
$ws.Range("M113").Value = -840.4484000000002
$ws.Range("N113").Value = -5950.3334
$ws.Range("H119").Value = 1595
$ws.Range("I119").Value = 465
$ws.Range("J119").Value = 3855
$ws.Range("K119").Value = 1395
$ws.Range("L119").Value = 11565
$ws.Range("M119").Value = 3443
$ws.Range("N119").Value = -21241

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 914.4545000000001
$ws.Range("I97").Value = 917.6667
$ws.Range("K97").Value = 917.6667
$ws.Range("M97").Value = -421.6667
$ws.Range("H98").Value = 34607.375
$ws.Range("J98").Value = 34607.375
$ws.Range("L98").Value = 34607.375
$ws.Range("N98").Value = -40597.375
$ws.Range("H126").Value = 1600
$ws.Range("I126").Value = 1600
$ws.Range("K126").Value = 4800
$ws.Range("M126").Value = -2330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 315.2143
$ws.Range("I46").Value = 299.18182
$ws.Range("K46").Value = 299.18182
$ws.Range("M46").Value = -111.18182
$ws.Range("H68").Value = 2940
$ws.Range("I68").Value = 2200
$ws.Range("K68").Value = 2200
$ws.Range("M68").Value = -1451
$ws.Range("H71").Value = 2940
$ws.Range("I71").Value = 2200
$ws.Range("K71").Value = 11000
$ws.Range("M71").Value = -7256

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 67500
$ws.Range("J109").Value = 67500
$ws.Range("L109").Value = 67500
$ws.Range("N109").Value = -70274
$ws.Range("H113").Value = 524.0769
$ws.Range("I113").Value = 401.375
$ws.Range("J113").Value = 720.4
$ws.Range("K113").Value = 1204.125
$ws.Range("L113").Value = 2161.2
$ws.Range("M113").Value = 965.875
$ws.Range("N113").Value = -6501.2
$ws.Range("H126").Value = 238937.17
$ws.Range("I126").Value = 435295.12
$ws.Range("J126").Value = 1240.6842
$ws.Range("K126").Value = 1305885.36
$ws.Range("L126").Value = 3722.0526
$ws.Range("M126").Value = -1303415.36
$ws.Range("N126").Value = -8662.052599999999
